$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I (I0) and J (IF), matching the style of the
# existing header row (copy formatting from H1 which already carries it).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I0 (col I) and IF (col J), rows 2-27
$i0 = @(8, 9, 9, 7, 8, 7, 7, 9, 7, 4, 9, 7, 7, 7, 6, 7, 1, 1, 9, 4, 8, 9, 8, 5, 7, 8)
$if = @(8, 9, 9, 7, 8, 7, 7, 9, 7, 5, 9, 7, 8, 8, 6, 8, 2, 2, 9, 5, 8, 9, 8, 6, 7, 8)

for ($r = 2; $r -le 27; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $i0[$idx]
    $ws.Cells.Item($r, 10).Value = $if[$idx]
}
